$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = -7.983999999999994
$ws.Range("C7").Value = -13.5504
$ws.Range("E7").Value = 15.5146
$ws.Range("B9").Value = 5.998099999999996
$ws.Range("E10").Value = 16.5827
$ws.Range("C12").Value = -10.6412
$ws.Range("B13").Value = 6.646799999999998
$ws.Range("E13").Value = 16.5236
$ws.Range("C14").Value = -13.42899999999999
$ws.Range("D15").Value = -8.752499999999996
$ws.Range("B16").Value = 5.226600000000002
$ws.Range("E16").Value = 15.8038
$ws.Range("B18").Value = 5.786499999999998
$ws.Range("C19").Value = -12.25420000000001
$ws.Range("B20").Value = 8.858399999999996
$ws.Range("E20").Value = 15.65419999999999
$ws.Range("E24").Value = 16.68070000000001
$ws.Range("B26").Value = 7.025900000000003
$ws.Range("C26").Value = -13.3047
$ws.Range("B27").Value = 6.5378
$ws.Range("C27").Value = -13.39329999999999
$ws.Range("D28").Value = -8.818299999999997
$ws.Range("B29").Value = 5.162199999999997
$ws.Range("C29").Value = -10.7708
$ws.Range("E32").Value = 16.84709999999998
$ws.Range("D33").Value = -7.759400000000004
$ws.Range("B35").Value = 8.551800000000004
$ws.Range("D35").Value = -8.878299999999991
$ws.Range("B36").Value = 9.186100000000009
$ws.Range("C37").Value = -13.2575
$ws.Range("C38").Value = -13.2749
$ws.Range("D38").Value = -8.666499999999996
$ws.Range("E39").Value = 16.171
$ws.Range("D43").Value = -8.253700000000002
$ws.Range("D44").Value = -7.632100000000001
$ws.Range("B45").Value = 5.609500000000001
$ws.Range("D45").Value = -7.990899999999995
$ws.Range("C47").Value = -11.57700000000001
$ws.Range("D47").Value = -7.461000000000001
$ws.Range("E47").Value = 16.9571
$ws.Range("E48").Value = 17.4634
$ws.Range("C51").Value = -12.6206
$ws.Range("D51").Value = -7.919800000000002
$ws.Range("C52").Value = -11.3441
$ws.Range("E52").Value = 17.32850000000001
$ws.Range("D54").Value = -8.322799999999999
$ws.Range("B55").Value = 5.806999999999996
$ws.Range("C55").Value = -14.04770000000001
$ws.Range("E56").Value = 16.41870000000001
$ws.Range("B57").Value = 5.052799999999998
$ws.Range("D57").Value = -8.433299999999997
$ws.Range("D62").Value = -8.400299999999994
$ws.Range("D63").Value = -7.887099999999997
$ws.Range("D67").Value = -6.545200000000003
$ws.Range("B69").Value = 6.184699999999993
$ws.Range("C69").Value = -11.2874
$ws.Range("C70").Value = -12.232
$ws.Range("D70").Value = -8.078900000000004
$ws.Range("B76").Value = 5.295600000000004
$ws.Range("C76").Value = -12.5997
$ws.Range("B78").Value = 9.7254
$ws.Range("C81").Value = -12.3951
$ws.Range("D81").Value = -8.258700000000005
$ws.Range("B82").Value = 5.8905
$ws.Range("B83").Value = 6.452599999999998
$ws.Range("C83").Value = -14.171
$ws.Range("E84").Value = 17.04459999999999
$ws.Range("D88").Value = -7.361599999999995
$ws.Range("B93").Value = 6.107199999999997
$ws.Range("C94").Value = -10.0283
$ws.Range("D96").Value = -8.268299999999996
$ws.Range("B97").Value = 6.075400000000001
$ws.Range("D99").Value = -7.934799999999997
$ws.Range("C100").Value = -10.9208
$ws.Range("E100").Value = 16.8979
$ws.Range("E101").Value = 16.65050000000003
$ws.Range("C102").Value = -13.76410000000001
